$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C15").Value = "Qtd Vr Irregular"
$ws.Range("G15").Value = "Qtd Vr Irregular"
$ws.Range("B15").Value = "Mês"
$ws.Range("F15").Value = "Mês"
$ws.Range("G6").Value = "MÍNIMO"
$ws.Range("H6").Value = "MÍNIMO"
$ws.Range("G8").Value = "MÁXIMO"
$ws.Range("H8").Value = "MÁXIMO"
